# ---------------------------------------------------------------------------
# Commit: "Thu, Mar 19, 2020  2:04:25 AM"
#
# The canonical diff swaps the content of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml (the deck's slide-master theme colors change from the
# "Integral"/Red-Violet palette to the built-in "Office" palette), and it
# also re-points a table on slide 5 from the deck's custom table style to
# the built-in "Medium Style 2 - Accent 1" table style
# ({3916D286-DCC7-4568-AB89-4E88283FECFD}).
#
# This host only exposes one live/editable theme (the one driving the slide
# master / the visible slides), reached via SlideMaster.Theme.ThemeColorScheme,
# so we rewrite its 12 theme colors to the "Office" palette that the diff
# shows ending up in the master's theme part.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

function HexToComRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("after") theme colors for the master's theme: built-in "Office"
# color scheme, in dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order.
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$masterTheme = $p.SlideMaster.Theme
$colorScheme = $masterTheme.ThemeColorScheme

for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $colorScheme.Item($i + 1).RGB = HexToComRGB $officeColors[$i]
}

# Slide 5's table (Google Shape;122;p17) switches from the deck's custom
# table style to the built-in "Medium Style 2 - Accent 1" style.
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{3916D286-DCC7-4568-AB89-4E88283FECFD}")
